$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run
#   "Possibility to access the different actions for an act by a slide left on the "
# into
#   "Possibility to access the"  +  " different actions for an act by a slide left on the "
# and move the (hidden) _GoBack bookmark so it still sits between the two
# halves (it previously sat right after the full run).
# ---------------------------------------------------------------------------

$needle = "Possibility to access the different actions for an act by a slide left on the "
$firstHalf = "Possibility to access the"

$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $runStart = $rng.Start
    $splitAt = $runStart + $firstHalf.Length

    # Force Word to split the run in two at $splitAt without altering the
    # visible character formatting: flip Bold on then back off over the
    # first half only.
    $firstPart = $d.Range($runStart, $splitAt)
    $firstPart.Bold = 1
    $firstPart.Bold = 0

    # Relocate the _GoBack bookmark to the new split point (it is hidden
    # from the Bookmarks enumeration but still reachable by name).
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
    $d.Bookmarks.Add("_GoBack", $d.Range($splitAt, $splitAt))
}

# ---------------------------------------------------------------------------
# Change 2: append a new trailing paragraph containing a single space.
# ---------------------------------------------------------------------------

$endRng = $d.Range($d.Content.End, $d.Content.End)
$endRng.InsertXML("<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`"><w:r><w:t> </w:t></w:r></w:p>")
